$p = $ppt.ActivePresentation

# --- Slide 20: "Compound Statement" -----------------------------------
# Fix the curly-quote typo: the phrase should read
#   braces "{" and "}"
# with a closing smart quote before "and" (was an opening quote twice).
$s20 = $p.Slides.Item(20)
$shContent20 = $s20.Shapes.Item(2)
$tr20 = $shContent20.TextFrame.TextRange
$para1 = $tr20.Paragraphs(1)
$run3 = $para1.Runs(3)
$run3.Text = [char]0x201D + " and " + [char]0x201C

# --- Slide 27: "Procedure Call Statement" ------------------------------
# Combine the grammar rule for procedureCallStmt onto a single paragraph
# (append the trailing  ";" .  to the first rule and drop the standalone
# continuation paragraph), and shrink the grammar text from 18.5pt to
# 17.5pt to make room.
$s27 = $p.Slides.Item(27)
$shContent27 = $s27.Shapes.Item(2)
$tr27 = $shContent27.TextFrame.TextRange

# Paragraph 7: procedureCallStmt = procId "(" [ actualParameters ] ")"
$para7 = $tr27.Paragraphs(7)
$lastRun7 = $para7.Characters(50, 6)
$lastRun7.Text = ' ] ")" ";" .'
$para7.Font.Size = 17.5

# Paragraph 8: standalone continuation  ";" .   -- remove it entirely
$para8 = $tr27.Paragraphs(8)
$para8.Delete()

# Paragraph 9 (now 8 after the delete above): actualParameters = expression { "," expression } .
$para8b = $tr27.Paragraphs(8)
$para8b.Font.Size = 17.5
